# Weekly price update: insert two new rows (Primera / Segunda) for the
# Femacal de La Calera - Piña sheet, for the new reporting date 2023-04-05
# (Excel serial 45021). This shifts every existing row from row 912 onward
# down by two rows (dimension grows from A1:T978 to A1:T980).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 912 (keeps formatting of the
# surrounding rows, e.g. the date-number-format on column D).
$ws.Range("A912:T913").Insert()

# New row 912 - "Primera" quality
$ws.Range("A912").Value = 3
$ws.Range("B912").Value = "Femacal de La Calera"
$ws.Range("C912").Value = "Coquimbo"
$ws.Range("D912").Value = 45021
$ws.Range("E912").Value = 5
$ws.Range("F912").Value = "Fruta"
$ws.Range("G912").Value = 100108
$ws.Range("H912").Value = "Tropicales y subtropicales"
$ws.Range("I912").Value = 100108005
$ws.Range("J912").Value = "Piña"
$ws.Range("K912").Value = "Caramelo"
$ws.Range("L912").Value = "Primera"
$ws.Range("M912").Value = 162
$ws.Range("N912").Value = 21000
$ws.Range("O912").Value = 21000
$ws.Range("P912").Value = 21000
$ws.Range("Q912").Value = "$/caja 12 unidades"
$ws.Range("R912").Value = "Ecuador"
$ws.Range("S912").Value = 1750
$ws.Range("T912").Value = 12

# New row 913 - "Segunda" quality
$ws.Range("A913").Value = 3
$ws.Range("B913").Value = "Femacal de La Calera"
$ws.Range("C913").Value = "Coquimbo"
$ws.Range("D913").Value = 45021
$ws.Range("E913").Value = 5
$ws.Range("F913").Value = "Fruta"
$ws.Range("G913").Value = 100108
$ws.Range("H913").Value = "Tropicales y subtropicales"
$ws.Range("I913").Value = 100108005
$ws.Range("J913").Value = "Piña"
$ws.Range("K913").Value = "Caramelo"
$ws.Range("L913").Value = "Segunda"
$ws.Range("M913").Value = 108
$ws.Range("N913").Value = 21000
$ws.Range("O913").Value = 21000
$ws.Range("P913").Value = 21000
$ws.Range("Q913").Value = "$/caja 14 unidades"
$ws.Range("R913").Value = "Ecuador"
$ws.Range("S913").Value = 1500
$ws.Range("T913").Value = 14
